$d = $word.ActiveDocument

# The commit renames the two inline logo pictures that live in the
# document's headers/footers (Pearson logo in the footers, BTec logo in
# the headers):
#   footers: PearsonLogo picture  image1.png -> image2.png
#   headers: BTec_Logo-Orange picture  image2.jpg -> image1.jpg
#
# InlineShape has no settable Name in the Word object model, so each
# picture is briefly promoted to a floating Shape (which does expose
# Name), renamed, then converted back to an inline picture so the
# <wp:inline> layout is preserved.

function Rename-LogoInRange($rng, $newName) {
    $count = $rng.InlineShapes.Count
    for ($k = 1; $k -le $count; $k++) {
        $ils = $rng.InlineShapes($k)
        $shp = $ils.ConvertToShape()
        $shp.Name = $newName
        [void]$shp.ConvertToInlineShape()
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections($si)

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            Rename-LogoInRange $hdr.Range "image1.jpg"
        }

        $ftr = $sec.Footers($i)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            Rename-LogoInRange $ftr.Range "image2.png"
        }
    }
}
